$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Zero-out the forecast values for specific categories/rows that the
# commit changed to 0 (ecom, cafes, phone/internet, ira transfer, and
# the single mira activities row).
$rowsToZero = @(17..31) + @(77..91) + @(197..211) + @(272..286) + @(337)
foreach ($r in $rowsToZero) {
    $ws.Cells.Item($r, 3).Value = 0
}

# Remove the "lower" / "upper" columns (D and E) entirely.
$ws.Range("D:E").Delete()
